# This script applies the "Fruta / hortaliza, semanal" update:
# A brand-new weekly record is inserted at row 7 (shifting the existing
# rows 7-67 down to rows 8-68), and the new row 7 is populated with the
# new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 7, shifting rows 7-67 down to 8-68.
$ws.Rows(7).Insert(-4121)

# Populate the newly inserted row 7 with the new weekly record.
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C7").Value = "Los Lagos"
$ws.Range("D7").Value = 44530
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 100112026
$ws.Range("G7").Value = "Haba"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 180
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 13000
$ws.Range("N7").Value = '$/saco 25 kilos'
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 520
$ws.Range("Q7").Value = 25
$ws.Range("R7").Value = "Hortaliza"
